$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data record was inserted at row 318, pushing the existing rows
# 318-408 down to 319-409 (the sheet grows from A1:R408 to A1:R409).
$ws.Rows.Item(318).Insert()

$ws.Range("A318").Value = 6
$ws.Range("B318").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C318").Value = "Metropolitana"
$ws.Range("D318").Value = 44642
$ws.Range("E318").Value = 13
$ws.Range("F318").Value = 100112043
$ws.Range("G318").Value = "Pepino ensalada"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 230
$ws.Range("K318").Value = 16000
$ws.Range("L318").Value = 16000
$ws.Range("M318").Value = 16000
$ws.Range("N318").Value = "`$/caja 70 unidades"
$ws.Range("O318").Value = "Provincia de Limarí"
$ws.Range("P318").Value = 229
$ws.Range("Q318").Value = 70
$ws.Range("R318").Value = "Hortaliza"
